$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Version value from 1.0.0 to 1.1.0 (row with "Version" label in column A)
$ws.Range("B3").Value = "1.1.0"

# Update Date value to new timestamp (row with "Date" label in column A)
$ws.Range("B8").Value = "2023-07-10T23:08:03+02:00"
